$d = $word.ActiveDocument

$enDash = [char]0x2013
$rsquo  = [char]0x2019

# The to-do list currently has (in order):
#   4) "When I add too many ... Fix this."
#   5) "Resize table font and column width for smaller screens – media query."
#
# After the edit it should read:
#   4) "Resize table font and column width for smaller screens – media query."
#   5) "  Figure out why white-space: nowrap is not working."
#
# Replace paragraph 5 first (while its old text is still unique), then
# paragraph 4, so neither Find lands on the wrong paragraph.

$oldP5 = "Resize table font and column width for smaller screens $enDash media query."
$newP5 = "  Figure out why white-space: nowrap is not working."
$rng5 = $d.Paragraphs.Item(5).Range
$rng5.Find.Execute($oldP5, $true, $false, $false, $false, $false, $true, 1, $false, $newP5, 2) | Out-Null

$oldP4 = "When I add too many to do list items the screen doesn${rsquo}t expand to fit the additional rows $enDash they get hidden by the footer bar.  Fix this."
$newP4 = "Resize table font and column width for smaller screens $enDash media query."
$rng4 = $d.Paragraphs.Item(4).Range
$rng4.Find.Execute($oldP4, $true, $false, $false, $false, $false, $true, 1, $false, $newP4, 2) | Out-Null
